$wb = $excel.ActiveWorkbook

# --- Sheet "Overview" ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"              # zh-cn status
$wsOverview.Range("F2").Value = "Ready for handoff"              # de-de status
$wsOverview.Range("G2").Value = "2016-08-21 00:46:36"            # Latest HO Xliff Generate Date
$wsOverview.Columns.Item(5).ColumnWidth = 16.333333333333332
$wsOverview.Columns.Item(6).ColumnWidth = 16.333333333333332

# --- Sheet "zh-cn" ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"                  # Status
$wsZhCn.Range("H2").Value = "2016-08-21 00:46:33"                # Latest Handoff Datetime
$wsZhCn.Columns.Item(3).ColumnWidth = 16.333333333333332

# --- Sheet "de-de" ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"                  # Status
$wsDeDe.Range("H2").Value = "2016-08-21 00:46:36"                # Latest Handoff Datetime
$wsDeDe.Columns.Item(3).ColumnWidth = 16.333333333333332
